$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.88
$ws.Range("H2").Value = 5
$ws.Range("J2").Value = 3.7
$ws.Range("K2").Value = 3.75
$ws.Range("W2").Value = 2.12

$ws.Range("G3").Value = 5.7
$ws.Range("M3").Value = 1.09
$ws.Range("X3").Value = 12
$ws.Range("AC3").Value = 8.199999999999999
$ws.Range("AI3").Value = 1000

$ws.Range("F4").Value = 5.6
$ws.Range("H4").Value = 1.84
$ws.Range("I4").Value = 1.86
$ws.Range("J4").Value = 3.55
$ws.Range("K4").Value = 3.6
$ws.Range("P4").Value = 1.64
$ws.Range("U4").Value = 1.72
$ws.Range("AC4").Value = 8.800000000000001
$ws.Range("AH4").Value = 29

$ws.Range("I5").Value = 1.84
$ws.Range("T5").Value = 1.97
$ws.Range("U5").Value = 1.97
$ws.Range("Y5").Value = 8.6

$ws.Range("I6").Value = 2.42

$ws.Range("P7").Value = 1.95
$ws.Range("T7").Value = 1.76
$ws.Range("AK7").Value = 34
$ws.Range("AN7").Value = 40

$ws.Range("F8").Value = 1.69
$ws.Range("G8").Value = 1.7
$ws.Range("H8").Value = 5.7
$ws.Range("I8").Value = 6.2
$ws.Range("J8").Value = 4
$ws.Range("K8").Value = 4.3
$ws.Range("P8").Value = 2.14
$ws.Range("R8").Value = 1.44
$ws.Range("U8").Value = 2.04
$ws.Range("Y8").Value = 26
$ws.Range("Z8").Value = 60
$ws.Range("AB8").Value = 12.5
$ws.Range("AC8").Value = 11.5
$ws.Range("AD8").Value = 27
$ws.Range("AF8").Value = 12.5
$ws.Range("AG8").Value = 13
$ws.Range("AH8").Value = 21
$ws.Range("AI8").Value = 85
$ws.Range("AJ8").Value = 19.5
$ws.Range("AK8").Value = 20
$ws.Range("AL8").Value = 38
$ws.Range("AM8").Value = 120
$ws.Range("AN8").Value = 10
$ws.Range("AO8").Value = 100
